$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# 1) Rows 227-236: fill in the newly-tracked estimate/actual hours (G,H)
#    and status (K) for "week 17" tasks that were previously left blank.
# -------------------------------------------------------------------------
$ws.Range("G227").Value = 4
$ws.Range("H227").Value = 7
$ws.Range("K227").Value = "Pending"

$ws.Range("G228").Value = 2
$ws.Range("H228").Value = 4
$ws.Range("K228").Value = "Complete"

$ws.Range("G229").Value = 3
$ws.Range("H229").Value = 5
$ws.Range("K229").Value = "Pending"

$ws.Range("G230").Value = 3
$ws.Range("H230").Value = 4
$ws.Range("K230").Value = "Complete"

$ws.Range("G231").Value = 3
$ws.Range("K231").Value = "Not implemented"

$ws.Range("G232").Value = 4

$ws.Range("G233").Value = 3
$ws.Range("H233").Value = 3
$ws.Range("K233").Value = "Complete"

$ws.Range("G234").Value = 3
$ws.Range("H234").Value = 2
$ws.Range("K234").Value = "Pending"

$ws.Range("G235").Value = 4
$ws.Range("H235").Value = 5
$ws.Range("K235").Value = "Pending"

$ws.Range("G236").Value = 4
$ws.Range("H236").Value = 3
$ws.Range("K236").Value = "Complete"

# -------------------------------------------------------------------------
# 2) Row 237 starts a new "week 18" block, formatted like the other
#    week-header rows (e.g. row 227: bold, red-filled B:K, white text).
#    Copy that formatting down first, then fill in the values.
# -------------------------------------------------------------------------
$ws.Range("B227:K227").Copy() | Out-Null
$ws.Range("B237:K237").PasteSpecial(-4122) | Out-Null

# -------------------------------------------------------------------------
# 3) Rows 238-257 are regular data rows for week 18. Stamp them with the
#    same formatting used by other regular rows that already carry a date
#    in F/I (row 17 is a template with styles B19/C20/D2/E2/F24/G20/H20/
#    I24/J20/K2, exactly what the target needs).
# -------------------------------------------------------------------------
$ws.Range("B17:K17").Copy() | Out-Null
$ws.Range("B238:K257").PasteSpecial(-4122) | Out-Null

# -------------------------------------------------------------------------
# 4) Rows 258-265 are new trailing blank rows (same look as the blank rows
#    238-241 had before this edit - plain, no date format on F/I).
# -------------------------------------------------------------------------
$ws.Range("B238:K241").Copy() | Out-Null
$ws.Range("B258:K265").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Give every newly-materialised row (242-265) the same 20.25pt height as
# the rest of the table.
$ws.Range("242:265").RowHeight = 20.25

# -------------------------------------------------------------------------
# 5) Fill in the values for rows 237-257.
# -------------------------------------------------------------------------
$ws.Range("B237").Value = 18
$ws.Range("C237").Value = 1
$ws.Range("D237").Value = "Modify Physical view"
$ws.Range("E237").Value = "Architect"
$ws.Range("F237").Value = 42888
$ws.Range("G237").Value = 3
$ws.Range("H237").Value = 4
$ws.Range("I237").Value = 42918
$ws.Range("J237").Value = "HienNguyen"
$ws.Range("K237").Value = "Complete"

$ws.Range("B238").Value = 18
$ws.Range("C238").Value = 2
$ws.Range("D238").Value = "Modify static view"
$ws.Range("E238").Value = "Architect"
$ws.Range("F238").Value = 42888
$ws.Range("G238").Value = 3
$ws.Range("H238").Value = 5
$ws.Range("I238").Value = 42949
$ws.Range("J238").Value = "HaiTran"

$ws.Range("B239").Value = 18
$ws.Range("C239").Value = 3
$ws.Range("D239").Value = "Modify Dynamic view"
$ws.Range("E239").Value = "Architect"
$ws.Range("F239").Value = 42918
$ws.Range("G239").Value = 5
$ws.Range("I239").Value = 42980
$ws.Range("J239").Value = "HaiTran"

$ws.Range("B240").Value = 18
$ws.Range("C240").Value = 4
$ws.Range("D240").Value = "Draw class diagram function ""Manage Account"""
$ws.Range("E240").Value = "Detailed Design"
$ws.Range("F240").Value = 42918
$ws.Range("I240").Value = 43010
$ws.Range("J240").Value = "TaiNguyen"

$ws.Range("B241").Value = 18
$ws.Range("C241").Value = 5
$ws.Range("D241").Value = "Draw class diagram function ""Manage News"""
$ws.Range("E241").Value = "Detailed Design"
$ws.Range("F241").Value = 42918
$ws.Range("I241").Value = 43010
$ws.Range("J241").Value = "TaiNguyen"

$ws.Range("B242").Value = 18
$ws.Range("C242").Value = 6
$ws.Range("D242").Value = "Draw class diagram function ""Manage FAQ"""
$ws.Range("E242").Value = "Detailed Design"
$ws.Range("F242").Value = 42918
$ws.Range("I242").Value = 43010
$ws.Range("J242").Value = "MinhDoan"

$ws.Range("B243").Value = 18
$ws.Range("C243").Value = 7
$ws.Range("D243").Value = "Draw class diagram function ""Manage Homepage"""
$ws.Range("E243").Value = "Detailed Design"
$ws.Range("F243").Value = 42918
$ws.Range("I243").Value = 43010
$ws.Range("J243").Value = "KhoiNguyen"

$ws.Range("B244").Value = 18
$ws.Range("C244").Value = 8
$ws.Range("D244").Value = "Draw class diagram function ""Manage Categories"""
$ws.Range("E244").Value = "Detailed Design"
$ws.Range("F244").Value = 42918
$ws.Range("I244").Value = 43010
$ws.Range("J244").Value = "TaiNguyen"

$ws.Range("B245").Value = 18
$ws.Range("C245").Value = 9
$ws.Range("D245").Value = "Draw class diagram function ""Manage Banners"""
$ws.Range("E245").Value = "Detailed Design"
$ws.Range("F245").Value = 42918
$ws.Range("I245").Value = 43010
$ws.Range("J245").Value = "HienNguyen"

$ws.Range("B246").Value = 18
$ws.Range("C246").Value = 10
$ws.Range("D246").Value = "Draw class diagram function ""Manage Popups"""
$ws.Range("E246").Value = "Detailed Design"
$ws.Range("F246").Value = 42918
$ws.Range("I246").Value = 43010
$ws.Range("J246").Value = "MinhDoan"

$ws.Range("B247").Value = 18
$ws.Range("C247").Value = 11
$ws.Range("D247").Value = "Draw class diagram function ""Manage Languages"""
$ws.Range("E247").Value = "Detailed Design"
$ws.Range("F247").Value = 42918
$ws.Range("I247").Value = 43010
$ws.Range("J247").Value = "TaiNguyen"

$ws.Range("B248").Value = 18
$ws.Range("C248").Value = 12
$ws.Range("D248").Value = "Draw class diagram function ""Manage Examination"""
$ws.Range("E248").Value = "Detailed Design"
$ws.Range("F248").Value = 42918
$ws.Range("I248").Value = 43010
$ws.Range("J248").Value = "XuanThaiHien"

$ws.Range("B249").Value = 18
$ws.Range("C249").Value = 13
$ws.Range("D249").Value = "Draw sequence diagram function ""Manage Account"""
$ws.Range("E249").Value = "Detailed Design"
$ws.Range("F249").Value = 42918
$ws.Range("I249").Value = 43010
$ws.Range("J249").Value = "HaiTran"

$ws.Range("B250").Value = 18
$ws.Range("C250").Value = 14
$ws.Range("D250").Value = "Draw sequence diagram function ""Manage News"""
$ws.Range("E250").Value = "Detailed Design"
$ws.Range("F250").Value = 42918
$ws.Range("I250").Value = 43010
$ws.Range("J250").Value = "KhoiNguyen"

$ws.Range("B251").Value = 18
$ws.Range("C251").Value = 15
$ws.Range("D251").Value = "Draw sequence diagram function ""Manage Categories"""
$ws.Range("E251").Value = "Detailed Design"
$ws.Range("F251").Value = 42918
$ws.Range("I251").Value = 43010
$ws.Range("J251").Value = "MinhDoan"

$ws.Range("B252").Value = 18
$ws.Range("C252").Value = 16
$ws.Range("D252").Value = "Draw sequence diagram function ""Manage Banners"""
$ws.Range("E252").Value = "Detailed Design"
$ws.Range("F252").Value = 42918
$ws.Range("I252").Value = 43010
$ws.Range("J252").Value = "HienNguyen"

$ws.Range("B253").Value = 18
$ws.Range("C253").Value = 17
$ws.Range("D253").Value = "Draw sequence diagram function ""Manage Popups"""
$ws.Range("E253").Value = "Detailed Design"
$ws.Range("F253").Value = 42918
$ws.Range("I253").Value = 43010
$ws.Range("J253").Value = "MinhDoan"

$ws.Range("B254").Value = 18
$ws.Range("C254").Value = 18
$ws.Range("D254").Value = "Draw sequence diagram function ""Manage Homepage"""
$ws.Range("E254").Value = "Detailed Design"
$ws.Range("F254").Value = 42918
$ws.Range("I254").Value = 43010
$ws.Range("J254").Value = "HaiTran"

$ws.Range("B255").Value = 18
$ws.Range("C255").Value = 19
$ws.Range("D255").Value = "Draw sequence diagram function ""Manage Examination"""
$ws.Range("E255").Value = "Detailed Design"
$ws.Range("F255").Value = 42918
$ws.Range("I255").Value = 43010
$ws.Range("J255").Value = "HienNguyen"

$ws.Range("B256").Value = 18
$ws.Range("C256").Value = 20
$ws.Range("D256").Value = "Draw sequence diagram function ""Manage Language"""
$ws.Range("E256").Value = "Detailed Design"
$ws.Range("F256").Value = 42918
$ws.Range("I256").Value = 43010
$ws.Range("J256").Value = "XuanThaiHien"

$ws.Range("B257").Value = 18
$ws.Range("C257").Value = 21
$ws.Range("D257").Value = "Draw sequence diagram function ""Manage FAQ"""
$ws.Range("E257").Value = "Detailed Design"
$ws.Range("F257").Value = 42918
$ws.Range("I257").Value = 43010
$ws.Range("J257").Value = "XuanThaiHien"

# -------------------------------------------------------------------------
# 6) Extend the K-column status dropdown validation down to the new last
#    row, and move the view / selection to where the new rows now are.
# -------------------------------------------------------------------------
$ws.Range("K2:K241").Validation.Delete()
$ws.Range("K2:K265").Validation.Add(3, 1, 1, "=sta")

$ws.Range("D245").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 228
$aw.ScrollColumn = 1
